# "Generate Report for Handback"
# Updates the localization-status workbook after a handback: the status
# columns move from "Ready for handoff" to "Handed back: in sync with en-US",
# and the per-language sheets get their Latest Target File / Latest Handback
# File / Latest Handback DateTime columns populated (with a hyperlink on the
# Latest Target File cell), plus a couple of columns get wider to fit the
# new content.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: status columns (zh-cn / de-de) now read the new status
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Overview column widths for zh-cn / de-de status columns got wider
$overview.Range("E1").ColumnWidth = 29.2
$overview.Range("F1").ColumnWidth = 29.2

# ---------------------------------------------------------------------
# Helper that fills in the handback columns (I/J/K) for a language sheet
# and rewires the A2/A3/I2/I3 hyperlinks in the expected order.
# ---------------------------------------------------------------------
function Update-LanguageSheet {
    param(
        [string]$SheetName,
        [string]$XliffFile,
        [string]$HandbackDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Widen the Status / Latest Target File / Latest Handback File columns
    $ws.Range("C1").ColumnWidth = 29.2
    $ws.Range("I1").ColumnWidth = 39.2
    $ws.Range("J1").ColumnWidth = 39.2

    # Populate Latest Target File (I), Latest Handback File (J) and
    # Latest Handback DateTime (K) for the two data rows.
    $ws.Range("I2").Value = "5e1a45c3-46ab-4f7b-83e1-c01e94c7b632.md"
    $ws.Range("J2").Value = $XliffFile
    $ws.Range("K2").Value = $HandbackDateTime

    $ws.Range("I3").Value = "5e1a45c3-46ab-4f7b-83e1-c01e94c7b632.md"
    $ws.Range("J3").Value = $XliffFile
    $ws.Range("K3").Value = $HandbackDateTime

    # Rebuild the hyperlinks so the ordering/ids come out as
    # A2, I2, A3, I3 (matching how the handback report links the new
    # "Latest Target File" cells to the same source document).
    $targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b95aafe58545626dec426720d85197e34be51892/e2e/5e1a45c3-46ab-4f7b-83e1-c01e94c7b632.md"
    $ffffUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b95aafe58545626dec426720d85197e34be51892/e2e/ffff9b9a05c8-316b-403c-a34e-92095ea42bef.md"

    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $targetUrl, "", "", "5e1a45c3-46ab-4f7b-83e1-c01e94c7b632.md")
    $ws.Hyperlinks.Add($ws.Range("I2"), $targetUrl, "", "", "5e1a45c3-46ab-4f7b-83e1-c01e94c7b632.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), $ffffUrl, "", "", "ffff9b9a05c8-316b-403c-a34e-92095ea42bef.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), $targetUrl, "", "", "5e1a45c3-46ab-4f7b-83e1-c01e94c7b632.md")
}

Update-LanguageSheet "zh-cn" "5e1a45c3-46ab-4f7b-83e1-c01e94c7b632.868128d0013e40b92172d4f8000cc902d6794320.zh-cn.xlf" "2016-10-24 09:45:36"
Update-LanguageSheet "de-de" "5e1a45c3-46ab-4f7b-83e1-c01e94c7b632.868128d0013e40b92172d4f8000cc902d6794320.de-de.xlf" "2016-10-24 09:45:52"
